$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 55
$ws.Range("B55").Value = 6830656
$ws.Range("F55").Value = "Mazatlan FC Women"
$ws.Range("G55").Value = "Club Necaxa Women"
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 2
$ws.Range("J55").Value = "A"
$ws.Range("K55").Value = 3.5
$ws.Range("L55").Value = 3.5
$ws.Range("M55").Value = 1.833
$ws.Range("N55").Value = 4.333
$ws.Range("O55").Value = 3.75
$ws.Range("P55").Value = 1.727
$ws.Range("Q55").Value = 0.75
$ws.Range("R55").Value = 1.85
$ws.Range("S55").Value = 1.95
$ws.Range("U55").Value = 1.975
$ws.Range("V55").Value = 1.825
$ws.Range("W55").Value = -1
$ws.Range("Y55").Value = 0.7270000000000001
$ws.Range("Z55").Value = -1
$ws.Range("AA55").Value = 0.95
$ws.Range("AC55").Value = 0.825

# Row 56
$ws.Range("B56").Value = 6830657
$ws.Range("F56").Value = "Chivas Guadalajara Women"
$ws.Range("G56").Value = "Monterrey Women"
$ws.Range("H56").Value = 1
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = "H"
$ws.Range("K56").Value = 2.375
$ws.Range("L56").Value = 3.4
$ws.Range("M56").Value = 2.5
$ws.Range("N56").Value = 2.15
$ws.Range("O56").Value = 3.5
$ws.Range("P56").Value = 3
$ws.Range("Q56").Value = -0.25
$ws.Range("R56").Value = 1.875
$ws.Range("S56").Value = 1.925
$ws.Range("U56").Value = 1.75
$ws.Range("V56").Value = 1.95
$ws.Range("W56").Value = 1.15
$ws.Range("Y56").Value = -1
$ws.Range("Z56").Value = 0.875
$ws.Range("AA56").Value = -1
$ws.Range("AC56").Value = 0.95

# Row 109
$ws.Range("B109").Value = 6830711
$ws.Range("F109").Value = "Atletico San Luis Women"
$ws.Range("G109").Value = "Tijuana Women"
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 2
$ws.Range("J109").Value = "A"
$ws.Range("K109").Value = 3
$ws.Range("L109").Value = 3.6
$ws.Range("M109").Value = 2
$ws.Range("N109").Value = 4
$ws.Range("O109").Value = 3.8
$ws.Range("P109").Value = 1.666
$ws.Range("Q109").Value = 0.75
$ws.Range("R109").Value = 1.925
$ws.Range("S109").Value = 1.875
$ws.Range("T109").Value = 3
$ws.Range("W109").Value = -1
$ws.Range("Y109").Value = 0.6659999999999999
$ws.Range("Z109").Value = -1
$ws.Range("AA109").Value = 0.875

# Row 110
$ws.Range("B110").Value = 6830712
$ws.Range("F110").Value = "Tigres UANL Women"
$ws.Range("G110").Value = "Unam Pumas Women"
$ws.Range("H110").Value = 3
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = "H"
$ws.Range("K110").Value = 1.181
$ws.Range("L110").Value = 6
$ws.Range("M110").Value = 10
$ws.Range("N110").Value = 1.1
$ws.Range("O110").Value = 9
$ws.Range("P110").Value = 19
$ws.Range("Q110").Value = -2.5
$ws.Range("R110").Value = 1.8
$ws.Range("S110").Value = 2
$ws.Range("T110").Value = 3.75
$ws.Range("W110").Value = 0.1000000000000001
$ws.Range("Y110").Value = -1
$ws.Range("Z110").Value = 0.8
$ws.Range("AA110").Value = -1

# Row 229
$ws.Range("B229").Value = 7645712
$ws.Range("F229").Value = "Chivas Guadalajara Women"
$ws.Range("G229").Value = "Santos Laguna Women"
$ws.Range("H229").Value = 10
$ws.Range("I229").Value = 2
$ws.Range("J229").Value = "H"
$ws.Range("K229").Value = 1.1
$ws.Range("L229").Value = 8
$ws.Range("M229").Value = 13
$ws.Range("N229").Value = 1.03
$ws.Range("O229").Value = 17
$ws.Range("P229").Value = 41
$ws.Range("Q229").Value = -3.75
$ws.Range("R229").Value = 1.775
$ws.Range("S229").Value = 1.925
$ws.Range("T229").Value = 4.75
$ws.Range("U229").Value = 1.9
$ws.Range("V229").Value = 1.9
$ws.Range("W229").Value = 0.03000000000000003
$ws.Range("Y229").Value = -1
$ws.Range("Z229").Value = 0.7749999999999999
$ws.Range("AA229").Value = -1
$ws.Range("AB229").Value = 0.8999999999999999

# Row 231
$ws.Range("B231").Value = 7645781
$ws.Range("F231").Value = "Atletico San Luis Women"
$ws.Range("G231").Value = "Cruz Azul Women"
$ws.Range("H231").Value = 3
$ws.Range("I231").Value = 5
$ws.Range("J231").Value = "A"
$ws.Range("K231").Value = 2.1
$ws.Range("L231").Value = 3.6
$ws.Range("M231").Value = 2.8
$ws.Range("N231").Value = 2.45
$ws.Range("O231").Value = 3.6
$ws.Range("P231").Value = 2.375
$ws.Range("Q231").Value = 0
$ws.Range("R231").Value = 1.95
$ws.Range("S231").Value = 1.85
$ws.Range("T231").Value = 2.75
$ws.Range("U231").Value = 1.75
$ws.Range("V231").Value = 1.95
$ws.Range("W231").Value = -1
$ws.Range("Y231").Value = 1.375
$ws.Range("Z231").Value = -1
$ws.Range("AA231").Value = 0.8500000000000001
$ws.Range("AB231").Value = 0.75

# Row 245
$ws.Range("B245").Value = 7645794
$ws.Range("F245").Value = "Toluca Women"
$ws.Range("G245").Value = "Club America Women"
$ws.Range("K245").Value = 6
$ws.Range("L245").Value = 5
$ws.Range("M245").Value = 1.333
$ws.Range("N245").Value = 5.75
$ws.Range("O245").Value = 5
$ws.Range("P245").Value = 1.363
$ws.Range("Q245").Value = 1.5
$ws.Range("R245").Value = 1.825
$ws.Range("S245").Value = 1.975
$ws.Range("T245").Value = 3.25
$ws.Range("U245").Value = 1.825
$ws.Range("V245").Value = 1.975
$ws.Range("W245").Value = 4.75
$ws.Range("Z245").Value = 0.825
$ws.Range("AB245").Value = -0.5
$ws.Range("AC245").Value = 0.4875

# Row 246
$ws.Range("B246").Value = 7645793
$ws.Range("F246").Value = "Queretaro Women"
$ws.Range("G246").Value = "Cruz Azul Women"
$ws.Range("K246").Value = 1.909
$ws.Range("L246").Value = 3.5
$ws.Range("M246").Value = 3.3
$ws.Range("N246").Value = 2.1
$ws.Range("O246").Value = 3.4
$ws.Range("P246").Value = 2.875
$ws.Range("Q246").Value = -0.25
$ws.Range("R246").Value = 1.875
$ws.Range("S246").Value = 1.925
$ws.Range("T246").Value = 2.75
$ws.Range("U246").Value = 2
$ws.Range("V246").Value = 1.8
$ws.Range("W246").Value = 1.1
$ws.Range("Z246").Value = 0.875
$ws.Range("AB246").Value = 0.5
$ws.Range("AC246").Value = -0.5

# Row 251
$ws.Range("B251").Value = 7926076
$ws.Range("F251").Value = "Leon Women"
$ws.Range("G251").Value = "Monterrey Women"
$ws.Range("H251").Value = 1
$ws.Range("I251").Value = 4
$ws.Range("K251").Value = 4.5
$ws.Range("L251").Value = 4
$ws.Range("M251").Value = 1.571
$ws.Range("N251").Value = 4.75
$ws.Range("O251").Value = 4
$ws.Range("P251").Value = 1.533
$ws.Range("Q251").Value = 1
$ws.Range("R251").Value = 1.825
$ws.Range("S251").Value = 1.975
$ws.Range("U251").Value = 1.975
$ws.Range("V251").Value = 1.825
$ws.Range("Y251").Value = 0.5329999999999999
$ws.Range("AA251").Value = 0.9750000000000001
$ws.Range("AB251").Value = 0.9750000000000001
$ws.Range("AC251").Value = -1

# Row 252
$ws.Range("B252").Value = 7645798
$ws.Range("F252").Value = "Atlas Women"
$ws.Range("G252").Value = "Unam Pumas Women"
$ws.Range("H252").Value = 0
$ws.Range("I252").Value = 2
$ws.Range("K252").Value = 2.4
$ws.Range("L252").Value = 3.6
$ws.Range("M252").Value = 2.4
$ws.Range("N252").Value = 2.375
$ws.Range("O252").Value = 3.75
$ws.Range("P252").Value = 2.375
$ws.Range("Q252").Value = 0
$ws.Range("R252").Value = 1.925
$ws.Range("S252").Value = 1.875
$ws.Range("U252").Value = 1.825
$ws.Range("V252").Value = 1.975
$ws.Range("Y252").Value = 1.375
$ws.Range("AA252").Value = 0.875
$ws.Range("AB252").Value = -1
$ws.Range("AC252").Value = 0.9750000000000001

# Row 278
$ws.Range("B278").Value = 7645820
$ws.Range("E278").Value = 45381.83333333334
$ws.Range("F278").Value = "Atletico San Luis Women"
$ws.Range("G278").Value = "Pachuca Women"
$ws.Range("K278").Value = 9
$ws.Range("L278").Value = 4
$ws.Range("M278").Value = 1.333
$ws.Range("N278").Value = 8.5
$ws.Range("O278").Value = 5
$ws.Range("P278").Value = 1.25
$ws.Range("Q278").Value = 1.75
$ws.Range("R278").Value = 1.95
$ws.Range("S278").Value = 1.85
$ws.Range("T278").Value = 3.5
$ws.Range("U278").Value = 1.85
$ws.Range("V278").Value = 1.95

# Row 279
$ws.Range("B279").Value = 8030105
$ws.Range("E279").Value = 45382.00694444445
$ws.Range("F279").Value = "Juarez FC Women"
$ws.Range("G279").Value = "Toluca Women"
$ws.Range("K279").Value = 1.8
$ws.Range("L279").Value = 3.75
$ws.Range("M279").Value = 3.75
$ws.Range("N279").Value = 1.8
$ws.Range("O279").Value = 3.75
$ws.Range("P279").Value = 3.75
$ws.Range("Q279").Value = -0.5
$ws.Range("R279").Value = 1.825
$ws.Range("S279").Value = 1.975

# Row 280
$ws.Range("B280").Value = 7645819
$ws.Range("E280").Value = 45382.625
$ws.Range("F280").Value = "Chivas Guadalajara Women"
$ws.Range("G280").Value = "Atlas Women"
$ws.Range("K280").Value = 1.222
$ws.Range("L280").Value = 6
$ws.Range("M280").Value = 8
$ws.Range("P280").Value = 9.5
$ws.Range("R280").Value = 1.975
$ws.Range("S280").Value = 1.825
$ws.Range("T280").Value = 3.25
$ws.Range("U280").Value = 1.925
$ws.Range("V280").Value = 1.875

# Row 281
$ws.Range("B281").Value = 7645730
$ws.Range("E281").Value = 45382.83333333334
$ws.Range("F281").Value = "Puebla Women"
$ws.Range("G281").Value = "Tigres UANL Women"
$ws.Range("K281").Value = 13
$ws.Range("L281").Value = 6.5
$ws.Range("M281").Value = 1.166
$ws.Range("N281").Value = 19
$ws.Range("O281").Value = 7.5
$ws.Range("P281").Value = 1.111
$ws.Range("Q281").Value = 2.5
$ws.Range("R281").Value = 1.8
$ws.Range("S281").Value = 2
$ws.Range("T281").Value = 3.5
$ws.Range("U281").Value = 1.9
$ws.Range("V281").Value = 1.9

# Add new row 282 (match between Santos Laguna Women and Unam Pumas Women)
$ws.Range("A281").Copy()
$ws.Range("A282").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("E281").Copy()
$ws.Range("E282").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A282").Value = 280
$ws.Range("B282").Value = 7645731
$ws.Range("C282").Value = "Mexico Liga MX Femenil"
$ws.Range("D282").Value = "Mexico Liga MX Femenil"
$ws.Range("E282").Value = 45382.92083333333
$ws.Range("F282").Value = "Santos Laguna Women"
$ws.Range("G282").Value = "Unam Pumas Women"
$ws.Range("K282").Value = 17
$ws.Range("L282").Value = 7
$ws.Range("M282").Value = 1.125
$ws.Range("N282").Value = 17
$ws.Range("O282").Value = 7
$ws.Range("P282").Value = 1.125
$ws.Range("Q282").Value = 2.25
$ws.Range("R282").Value = 1.9
$ws.Range("S282").Value = 1.9
$ws.Range("T282").Value = 3.25
$ws.Range("U282").Value = 1.9
$ws.Range("V282").Value = 1.9
$ws.Range("W282").Value = 0
$ws.Range("X282").Value = 0
$ws.Range("Y282").Value = 0
$ws.Range("Z282").Value = 0
$ws.Range("AA282").Value = 0
